$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $escaped = $text -replace '"', '""'
    $r.Formula = '="' + $escaped + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

Set-TextValue 'D2' '26.457.34'
Set-TextValue 'E2' '  +3.48%  '
Set-TextValue 'D3' '1.729.20'
Set-TextValue 'E3' '  +3.83%  '
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  +0.23%  '
Set-TextValue 'D5' '243.53'
Set-TextValue 'E5' '  +2.62%  '
Set-TextValue 'D6' '1.001'
Set-TextValue 'E6' '  +0.09%  '
Set-TextValue 'D7' '0.4790'
Set-TextValue 'E7' '  +3.66%  '
Set-TextValue 'D8' '0.2661'
Set-TextValue 'E8' '  +3.15%  '
Set-TextValue 'D9' '0.06218'
Set-TextValue 'E9' '  +1.26%  '
Set-TextValue 'D10' '1.733.89'
Set-TextValue 'E10' '  +4.15%  '
Set-TextValue 'D11' '0.07126'
Set-TextValue 'E11' '  +2.87%  '
Set-TextValue 'D12' '15.69'
Set-TextValue 'E12' '  +5.91%  '
Set-TextValue 'D13' '0.6140'
Set-TextValue 'E13' '  +7.25%  '
Set-TextValue 'D14' '4.529'
Set-TextValue 'E14' '  +4.25%  '
Set-TextValue 'D15' '76.86'
Set-TextValue 'E15' '  +2.23%  '
Set-TextValue 'E16' '  +0.08%  '
Set-TextValue 'D17' '26.475.67'
Set-TextValue 'E17' '  +3.56%  '
Set-TextValue 'D18' '1.001'
Set-TextValue 'E18' '  +0.15%  '
Set-TextValue 'D19' '0.000006899'
Set-TextValue 'E19' '  +3.21%  '
Set-TextValue 'D20' '11.70'
Set-TextValue 'E20' '  +2.81%  '
Set-TextValue 'D21' '1.956.79'
Set-TextValue 'E21' '  +4.32%  '
Set-TextValue 'D22' '4.570'
Set-TextValue 'E22' '  +3.15%  '
Set-TextValue 'D23' '8.893'
Set-TextValue 'E23' '  +3.33%  '
Set-TextValue 'D24' '5.318'
Set-TextValue 'E24' '  +1.91%  '
Set-TextValue 'D25' '136.23'
Set-TextValue 'E25' '  +1.50%  '
Set-TextValue 'E26' '  +2.81%  '
Set-TextValue 'D27' '1.791'
Set-TextValue 'E27' '  +4.09%  '
Set-TextValue 'D28' '1.400'
Set-TextValue 'E28' '  +1.80%  '
Set-TextValue 'D29' '106.73'
Set-TextValue 'E29' '  +2.46%  '
Set-TextValue 'D30' '3.978'
Set-TextValue 'E30' '  +0.93%  '
Set-TextValue 'D31' '0.08017'
Set-TextValue 'E31' '  +4.66%  '
Set-TextValue 'D32' '3.712'
Set-TextValue 'E32' '  +3.31%  '
Set-TextValue 'D33' '0.04536'
Set-TextValue 'E33' '  +4.58%  '
Set-TextValue 'E34' '  +0.43%  '
Set-TextValue 'D35' '0.6350'
Set-TextValue 'E35' '  +4.78%  '
Set-TextValue 'D36' '0.9893'
Set-TextValue 'E36' '  +5.34%  '
Set-TextValue 'D37' '0.9326'
Set-TextValue 'E37' '  +0.73%  '
Set-TextValue 'D38' '109.64'
Set-TextValue 'E38' '  +2.20%  '
Set-TextValue 'D39' '1.975'
Set-TextValue 'E39' '  +7.76%  '
Set-TextValue 'D40' '2.417'
Set-TextValue 'E40' '  -0.75%  '
Set-TextValue 'E41' '  +0.62%  '
Set-TextValue 'D42' '0.01505'
Set-TextValue 'E42' '  +3.85%  '
Set-TextValue 'D43' '5.651'
Set-TextValue 'E43' '  +11.70%  '
Set-TextValue 'D44' '0.3896'
Set-TextValue 'E44' '  +5.05%  '
Set-TextValue 'D45' '6.944'
Set-TextValue 'E45' '  +13.88%  '
Set-TextValue 'D46' '0.1188'
Set-TextValue 'E46' '  +7.07%  '
Set-TextValue 'D47' '0.05331'
Set-TextValue 'E47' '  +1.23%  '
Set-TextValue 'D48' '7.884'
Set-TextValue 'E48' '  +4.10%  '
Set-TextValue 'D49' '30.76'
Set-TextValue 'E49' '  -0.81%  '
Set-TextValue 'D50' '1.268'
Set-TextValue 'E50' '  +5.11%  '
Set-TextValue 'D51' '0.3418'
Set-TextValue 'E51' '  +3.02%  '
